$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.1467555
$ws.Range("H2").Value = 20.293511
$ws.Range("I2").Value = 0.1516003594919049
$ws.Range("J2").Value = 0.1102643619993968
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.032708
$ws.Range("N2").Value = 0.065416
$ws.Range("Q2").Value = 0.3318800788940001
$ws.Range("R2").Value = 1.327520315576
$ws.Range("S2").Value = 0.1516003594919049
$ws.Range("T2").Value = 0.1102643619993968

$ws.Range("I3").Value = 0.1339722830802056
$ws.Range("J3").Value = 0.1461642475877201
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.032708
$ws.Range("N3").Value = 0.065416
$ws.Range("Q3").Value = 0.2932890926333334
$ws.Range("R3").Value = 1.7597345558
$ws.Range("S3").Value = 0.1339722830802056
$ws.Range("T3").Value = 0.1461642475877201

$ws.Range("G4").Value = 17.36323866666666
$ws.Range("H4").Value = 52.089716
$ws.Range("I4").Value = 0.2594201884346587
$ws.Range("J4").Value = 0.2830283681059314
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.032708
$ws.Range("N4").Value = 0.065416
$ws.Range("Q4").Value = 0.5679168103093333
$ws.Range("R4").Value = 3.407500861856
$ws.Range("S4").Value = 0.2594201884346587
$ws.Range("T4").Value = 0.2830283681059314

$ws.Range("G5").Value = 6.601931
$ws.Range("H5").Value = 13.203862
$ws.Range("I5").Value = 0.09863794519743292
$ws.Range("J5").Value = 0.07174290438742112
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.032708
$ws.Range("N5").Value = 0.065416
$ws.Range("Q5").Value = 0.215935959148
$ws.Range("R5").Value = 0.863743836592
$ws.Range("S5").Value = 0.09863794519743292
$ws.Range("T5").Value = 0.07174290438742112

$ws.Range("G6").Value = 5.391932333333334
$ws.Range("H6").Value = 16.175797
$ws.Range("I6").Value = 0.080559631114533
$ws.Range("J6").Value = 0.0878908502346763
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.032708
$ws.Range("N6").Value = 0.065416
$ws.Range("Q6").Value = 0.1763593227586667
$ws.Range("R6").Value = 1.058155936552
$ws.Range("S6").Value = 0.080559631114533
$ws.Range("T6").Value = 0.0878908502346763

$ws.Range("G7").Value = 18.460197
$ws.Range("H7").Value = 55.380591
$ws.Range("I7").Value = 0.2758095926812649
$ws.Range("J7").Value = 0.3009092676848542
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.032708
$ws.Range("N7").Value = 0.065416
$ws.Range("Q7").Value = 0.6037961234760001
$ws.Range("R7").Value = 3.622776740856
$ws.Range("S7").Value = 0.2758095926812649
$ws.Range("T7").Value = 0.3009092676848542
